$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '72.356.16'
$ws.Range("E2").Value = '  +4.40%  '

$ws.Range("D3").Value = '3.632.38'
$ws.Range("E3").Value = '  +7.16%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").Value = "'597.52"
$ws.Range("E5").Value = '  +1.93%  '

$ws.Range("D6").Value = "'182.49"
$ws.Range("E6").Value = '  +0.99%  '

$ws.Range("D7").Value = '3.622.38'
$ws.Range("E7").Value = '  +7.19%  '

$ws.Range("D8").Value = "'0.607"
$ws.Range("E8").Value = '  +1.76%  '

$ws.Range("E9").Value = '  +0.12%  '

$ws.Range("E10").Value = '  +4.90%  '

$ws.Range("D11").Value = "'0.610"
$ws.Range("E11").Value = '  +2.73%  '

$ws.Range("D12").Value = "'50.43"
$ws.Range("E12").Value = '  +3.66%  '

$ws.Range("D13").Value = "'0.0000291"
$ws.Range("E13").Value = '  +2.34%  '

$ws.Range("D14").Value = "'706.46"
$ws.Range("E14").Value = '  +3.95%  '

$ws.Range("D15").Value = '4.210.63'
$ws.Range("E15").Value = '  +7.06%  '

$ws.Range("E16").Value = '  +3.70%  '

$ws.Range("D17").Value = '72.503.24'
$ws.Range("E17").Value = '  +4.54%  '

$ws.Range("D18").Value = '3.564.96'
$ws.Range("E18").Value = '  +4.72%  '

$ws.Range("E19").Value = '  +1.99%  '

$ws.Range("D20").Value = "'18.61"

$ws.Range("E21").Value = '  +3.77%  '

$ws.Range("D22").Value = "'0.935"
$ws.Range("E22").Value = '  +3.13%  '

$ws.Range("D23").Value = "'5.85"
$ws.Range("E23").Value = '  +7.86%  '

$ws.Range("E24").Value = '  +3.88%  '

$ws.Range("D25").Value = "'104.76"
$ws.Range("E25").Value = '  +1.65%  '

$ws.Range("E26").Value = '  +2.91%  '

$ws.Range("E27").Value = '  +5.11%  '

$ws.Range("D28").Value = "'10.04"
$ws.Range("E28").Value = '  +4.18%  '

$ws.Range("D29").Value = "'35.49"
$ws.Range("E29").Value = '  +4.67%  '

$ws.Range("E30").Value = '  +3.98%  '

$ws.Range("D31").Value = "'7.46"
$ws.Range("E31").Value = '  +7.09%  '

$ws.Range("D32").Value = "'4.18"
$ws.Range("E32").Value = '  +15.72%  '

$ws.Range("D33").Value = "'595.93"
$ws.Range("E33").Value = '  +7.11%  '

$ws.Range("E34").Value = '  +1.80%  '

$ws.Range("D35").Value = "'0.109"
$ws.Range("E35").Value = '  +2.17%  '

$ws.Range("D36").Value = "'59.83"
$ws.Range("E36").Value = '  +2.18%  '

$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = '  -0.05%  '

$ws.Range("E38").Value = '  +4.02%  '

$ws.Range("D39").Value = '3.640.07'
$ws.Range("E39").Value = '  -0.83%  '

$ws.Range("E40").Value = '  +8.13%  '

$ws.Range("D41").Value = "'35.92"
$ws.Range("E41").Value = '  +0.94%  '

$ws.Range("E42").Value = '  +6.76%  '

$ws.Range("E43").Value = '  +5.28%  '

$ws.Range("D44").Value = "'0.0451"
$ws.Range("E44").Value = '  +6.53%  '

$ws.Range("E45").Value = '  +3.42%  '

$ws.Range("E46").Value = '  +3.85%  '

$ws.Range("E47").Value = '  +5.18%  '

$ws.Range("E48").Value = '  +5.38%  '

$ws.Range("E49").Value = '  +1.64%  '

$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = '  -0.10%  '

$ws.Range("D51").Value = "'133.88"
$ws.Range("E51").Value = '  +0.18%  '
